$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: "4-1" / "4-2" / "4-3" entered as quote-prefixed (text) values.
# A4 additionally carries a date-like number format (d-mmm) left over from
# Excel's auto-detection of the "4-1" pattern, combined with the quote
# prefix that forces it to stay text.
$ws.Range("A4").Value = "'4-1"
$ws.Range("A4").NumberFormat = "d-mmm"
$ws.Range("B4").Value = "'4-2"
$ws.Range("C4").Value = "'4-3"

# Row 5: example of looking up a matched value's adjacent cells -
# username/password entered first, then the matched "Firstname Lastname".
$ws.Range("B5").Value = "username"
$ws.Range("C5").Value = "password"
$ws.Range("A5").Value = "Firstname Lastname"

# Widen column A (names) and column C (passwords) so the new text fits.
$ws.Columns.Item(1).ColumnWidth = 20.6666666667
$ws.Columns.Item(3).ColumnWidth = 13.3333333333

# Move the active selection below the new data, as in the saved workbook.
$null = $ws.Range("A6").Select()
